$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 4th data row with client "Quam Ullamcorper"
$ws.Range("A4").Value = "Quam Ullamcorper"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "Justo@Lorem.com"
$ws.Range("D4").Value = 61231212
$ws.Range("E4").Value = "Quam Adipiscing Parturient Justo 123"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "05/06/2014 - 18:25:52"

# Resize columns: A wider (for longer name), E wider (for longer address)
# (COM ColumnWidth is quantized to whole pixels at the engine's Maximum Digit
# Width of 7, same as real Excel; the literal input values below are the
# ones whose resulting (quantized) stored column width lands closest to the
# target stored widths of 12.989887640449439 and 26.18988764044944.)
$ws.Columns.Item(1).ColumnWidth = 12.29
$ws.Columns.Item(5).ColumnWidth = 25.4
